$d = $word.ActiveDocument

# Locate the original sentence that needs to be restructured into multiple runs.
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute(", has been developed by respected Arduino and Processing community members seeking to use the hardware in a configuration similar to that in our project.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate target sentence to edit."
}

$start = $rng.Start

# Remove the original text entirely; we will rebuild it as several runs.
$rng.Text = ""

# The replacement text, split into the pieces that must become separate <w:r> runs.
$parts = @(
    ", ",
    "was",
    " developed by Arduino and Processing community members ",
    "for a hardware",
    " configuration similar to that ",
    "of",
    " our project."
)

# Insert each piece after the previous one, advancing the insertion point each time.
# Because the (soon to be relocated) _GoBack bookmark still sits at $start, each
# InsertAfter call is forced to create a brand-new run instead of being merged into
# the preceding text.
$cur = $d.Range($start, $start)
foreach ($part in $parts) {
    $cur.InsertAfter($part)
    $cur = $d.Range($cur.End, $cur.End)
}

$realEnd = $cur.End

# Temporarily append a placeholder character after the real content so that $realEnd
# is no longer the very last character position of the paragraph (that special
# position confuses Bookmarks.Add). Then move the _GoBack bookmark to $realEnd and
# delete the placeholder, leaving the bookmark correctly collapsed right after the
# newly inserted text.
$cur.InsertAfter("Z")

$bm = $d.Bookmarks("_GoBack")
$bm.Delete()
$d.Bookmarks.Add("_GoBack", $d.Range($realEnd, $realEnd))

$placeholderRng = $d.Range($realEnd, $realEnd + 1)
$placeholderRng.Text = ""
